# "Desplazamiento de todas las filas N posiciones desde un punto dado"
#
# 1. Rename the original sheet "Hoja1" -> "Template Tabla".
# 2. Insert a new sheet "Desplazar Filas" right after it, and fill it with
#    the list of months used by the row-shifting helper (A1:A10).
# 3. Move the selection on the template sheet to B7 (where the new
#    "shift rows" entry point lives) and leave the template sheet active.

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename ---------------------------------------------------
$wsTemplate = $wb.Worksheets.Item(1)
$wsTemplate.Name = "Template Tabla"

# --- New sheet: "Desplazar Filas" --------------------------------------
$wsShift = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTemplate)
$wsShift.Name = "Desplazar Filas"

$meses = @("ENERO", "FEBRERO", "MARZO", "JUNIO", "JULIO", "AGOSTO", "SEPTIEMBRE", "OCTUBRE", "NOVIEMBRE", "DICIEMBRE")
for ($i = 0; $i -lt $meses.Length; $i++) {
    $wsShift.Cells.Item($i + 1, 1).Value = $meses[$i]
}

[void]$wsShift.Range("D26").Select()

# --- Back to the template sheet, park selection on B7 ------------------
$wsTemplate.Activate()
[void]$wsTemplate.Range("B7").Select()
